$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
try {
  $ws.Range("N500").Text = "5"
  Write-Output "text-ok"
} catch {
  Write-Output "text-err: $_"
}
